# Userstories workbook update: add "Loading Screen" user story, adjust
# row heights / view state to match the author's edit in Excel.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Tabelle1")
$ws2 = $wb.Worksheets.Item("Gamboard")

# --- New user story row on Tabelle1 -----------------------------------
$ws1.Range("A18").Value = "Loading Screen"
$ws1.Range("B18").Value = "Alexandra"
$ws1.Range("C18").Value = "Beim Start des Spiels sollte ein Ladescreen erscheinen der mich (zB mittels progressbar) darüber informiert wie weit das Spiel bereits geladen ist."
$ws1.Range("D18").Value = "Hoch"
$ws1.Rows.Item(18).RowHeight = 45

# --- Row height tweaks on existing rows (re-wrap adjustments) ---------
$ws1.Rows.Item(2).RowHeight = 30
$ws1.Rows.Item(6).RowHeight = 45

# --- Page setup on Tabelle1 --------------------------------------------
$ws1.PageSetup.Orientation = 1

# --- View / selection state --------------------------------------------
# Author ended up on Tabelle1, with C19 selected, so it becomes the
# active/selected tab (Gamboard loses tabSelected).
$ws2.Range("C11").Select()
$ws1.Activate()
$ws1.Range("C19").Select()
